$wb = $excel.ActiveWorkbook

# Sheet "展览" (first sheet) - update "想去人数" (F column) values
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 12499
$ws1.Range("F3").Value = 597
$ws1.Range("F4").Value = 2033
$ws1.Range("F5").Value = 263
$ws1.Range("F6").Value = 383
$ws1.Range("F8").Value = 12471
$ws1.Range("F9").Value = 3065
$ws1.Range("F10").Value = 530
$ws1.Range("F12").Value = 17
$ws1.Range("F14").Value = 21
$ws1.Range("F15").Value = 127
$ws1.Range("F16").Value = 643
$ws1.Range("F17").Value = 2834
$ws1.Range("F18").Value = 6068
$ws1.Range("F20").Value = 3598
$ws1.Range("F22").Value = 36

# Sheet "全部类型" (fourth sheet) - update "想去人数" (F column) values
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 12499
$ws4.Range("F3").Value = 597
$ws4.Range("F4").Value = 2033
$ws4.Range("F5").Value = 263
$ws4.Range("F7").Value = 383
$ws4.Range("F9").Value = 12471
$ws4.Range("F10").Value = 3065
$ws4.Range("F11").Value = 530
$ws4.Range("F13").Value = 17
$ws4.Range("F15").Value = 21
$ws4.Range("F16").Value = 127
$ws4.Range("F17").Value = 643
$ws4.Range("F18").Value = 2834
$ws4.Range("F20").Value = 6068
$ws4.Range("F22").Value = 3598
$ws4.Range("F24").Value = 36
